{"js": "const pairs = [\n  [\"2023-05-25 Thursday\", \"2023-05-26 Friday\"],\n  [\"77+21=98\", \"67+2=69\"],\n  [\"26-2=24\", \"44-0=44\"],\n  [\"0+82=82\", \"16+34=50\"],\n  [\"52+28=80\", \"83-42=41\"],\n  [\"45+54=99\", \"45-25=20\"],\n  [\"1+41=42\", \"3+32=35\"],\n  [\"86-1=85\", \"19+17=36\"],\n  [\"75-2=73\", \"42+13=55\"],\n  [\"22+6=28\", \"51+0=51\"],\n  [\"41+16=57\", \"47+37=84\"],\n  [\"44+26=70\", \"17+41=58\"],\n  [\"30-10=20\", \"23-20=3\"],\n  [\"78+13=91\", \"74-67=7\"],\n  [\"93-40=53\", \"36+6=42\"],\n  [\"87-82=5\", \"24-20=4\"],\n  [\"55-31=24\", \"20+11=31\"],\n  [\"49+39=88\", \"11-1=10\"],\n  [\"35+64=99\", \"72-22=50\"],\n  [\"44+36=80\", \"35+10=45\"],\n  [\"70-17=53\", \"91+7=98\"],\n  [\"10-2=8\", \"33+22=55\"],\n  [\"44+30=74\", \"19+72=91\"],\n  [\"73-43=30\", \"85-29=56\"],\n  [\"68-2=66\", \"83-24=59\"],\n  [\"42+44=86\", \"92-53=39\"],\n  [\"77+2=79\", \"8+44=52\"],\n  [\"25+6=31\", \"5+35=40\"],\n  [\"65-44=21\", \"89-15=74\"],\n  [\"36+38=74\", \"21+40=61\"],\n  [\"80+9=89\", \"45-32=13\"],\n  [\"89-80=9\", \"23+55=78\"],\n  [\"3+53=56\", \"45+40=85\"],\n  [\"27+42=69\", \"15+27=42\"],\n  [\"54+33=87\", \"38-12=26\"],\n  [\"61-20=41\", \"35+12=47\"],\n  [\"61-59=2\", \"13+36=49\"],\n  [\"45-22=23\", \"90-77=13\"],\n  [\"1+73=74\", \"1+19=20\"],\n  [\"52-28=24\", \"32+13=45\"],\n  [\"7+2=9\", \"88-78=10\"],\n  [\"41+56=97\", \"56-21=35\"],\n  [\"12+15=27\", \"14+15=29\"],\n  [\"99-42=57\", \"98-77=21\"],\n  [\"53+34=87\", \"25-3=22\"],\n  [\"38+32=70\", \"29-18=11\"],\n  [\"28+17=45\", \"6+82=88\"],\n  [\"57-15=42\", \"9+1=10\"],\n  [\"91-38=53\", \"19+49=68\"],\n  [\"49-27=22\", \"85-57=28\"],\n  [\"99-29=70\", \"54-54=0\"],\n  [\"53-13=40\", \"54-54=0\"],\n  [\"45-20=25\", \"16+20=36\"],\n  [\"55+6=61\", \"20+34=54\"],\n  [\"1+39=40\", \"83-34=49\"],\n  [\"73+5=78\", \"24-15=9\"],\n  [\"88+4=92\", \"17+56=73\"],\n  [\"95-7=88\", \"70-31=39\"],\n  [\"14+42=56\", \"36+14=50\"],\n  [\"14+6=20\", \"25+73=98\"],\n  [\"32-23=9\", \"14+40=54\"],\n  [\"60-6=54\", \"38+46=84\"],\n  [\"63+6=69\", \"62+35=97\"],\n  [\"60-7=53\", \"55-37=18\"],\n  [\"42+22=64\", \"49+10=59\"],\n  [\"43-18=25\", \"28+26=54\"],\n  [\"34+3=37\", \"0+73=73\"],\n  [\"56+10=66\", \"69-57=12\"],\n  [\"42+19=61\", \"41+51=92\"],\n  [\"56-24=32\", \"19+76=95\"],\n  [\"20-17=3\", \"49+0=49\"],\n  [\"77+12=89\", \"25+64=89\"],\n  [\"1+59=60\", \"99-54=45\"],\n  [\"19+47=66\", \"36-30=6\"],\n  [\"83-73=10\", \"28+44=72\"],\n  [\"35-5=30\", \"88-6=82\"],\n  [\"19+5=24\", \"96-33=63\"],\n  [\"62-17=45\", \"34-4=30\"],\n  [\"40-8=32\", \"3+66=69\"],\n  [\"33+57=90\", \"93-79=14\"],\n  [\"21-5=16\", \"93-55=38\"],\n  [\"61-52=9\", \"37+0=37\"],\n  [\"91-39=52\", \"53-41=12\"],\n  [\"67+8=75\", \"15+58=73\"],\n  [\"44+45=89\", \"86-0=86\"],\n  [\"18+14=32\", \"67-22=45\"],\n  [\"21+46=67\", \"46+27=73\"],\n  [\"20+25=45\", \"18+67=85\"],\n  [\"37-21=16\", \"66-39=27\"],\n  [\"45-44=1\", \"85-13=72\"],\n  [\"15-3=12\", \"32+58=90\"],\n  [\"28-27=1\", \"15+65=80\"],\n  [\"11+34=45\", \"11+51=62\"],\n  [\"69-42=27\", \"20-9=11\"],\n  [\"15+59=74\", \"62-32=30\"],\n  [\"50-3=47\", \"24+44=68\"],\n  [\"37+23=60\", \"75-46=29\"],\n  [\"13+8=21\", \"73-25=48\"],\n  [\"44+31=75\", \"47-0=47\"],\n  [\"33+12=45\", \"16+43=59\"],\n  [\"58-40=18\", \"36+18=54\"],\n];\n\nfor (const [before, after] of pairs) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n# Ordered list of (old, new) text replacements taken from the edit diff.\n# Every 'old' string is unique in the document, so a literal, non-wildcard\n# Find/Replace for each pair reproduces the target edit exactly.\n$pairs = @(\n  @('2023-05-25 Thursday', '2023-05-26 Friday'),\n  @('77+21=98', '67+2=69'),\n  @('26-2=24', '44-0=44'),\n  @('0+82=82', '16+34=50'),\n  @('52+28=80', '83-42=41'),\n  @('45+54=99', '45-25=20'),\n  @('1+41=42', '3+32=35'),\n  @('86-1=85', '19+17=36'),\n  @('75-2=73', '42+13=55'),\n  @('22+6=28', '51+0=51'),\n  @('41+16=57', '47+37=84'),\n  @('44+26=70', '17+41=58'),\n  @('30-10=20', '23-20=3'),\n  @('78+13=91', '74-67=7'),\n  @('93-40=53', '36+6=42'),\n  @('87-82=5', '24-20=4'),\n  @('55-31=24', '20+11=31'),\n  @('49+39=88', '11-1=10'),\n  @('35+64=99', '72-22=50'),\n  @('44+36=80', '35+10=45'),\n  @('70-17=53', '91+7=98'),\n  @('10-2=8', '33+22=55'),\n  @('44+30=74', '19+72=91'),\n  @('73-43=30', '85-29=56'),\n  @('68-2=66', '83-24=59'),\n  @('42+44=86', '92-53=39'),\n  @('77+2=79', '8+44=52'),\n  @('25+6=31', '5+35=40'),\n  @('65-44=21', '89-15=74'),\n  @('36+38=74', '21+40=61'),\n  @('80+9=89', '45-32=13'),\n  @('89-80=9', '23+55=78'),\n  @('3+53=56', '45+40=85'),\n  @('27+42=69', '15+27=42'),\n  @('54+33=87', '38-12=26'),\n  @('61-20=41', '35+12=47'),\n  @('61-59=2', '13+36=49'),\n  @('45-22=23', '90-77=13'),\n  @('1+73=74', '1+19=20'),\n  @('52-28=24', '32+13=45'),\n  @('7+2=9', '88-78=10'),\n  @('41+56=97', '56-21=35'),\n  @('12+15=27', '14+15=29'),\n  @('99-42=57', '98-77=21'),\n  @('53+34=87', '25-3=22'),\n  @('38+32=70', '29-18=11'),\n  @('28+17=45', '6+82=88'),\n  @('57-15=42', '9+1=10'),\n  @('91-38=53', '19+49=68'),\n  @('49-27=22', '85-57=28'),\n  @('99-29=70', '54-54=0'),\n  @('53-13=40', '54-54=0'),\n  @('45-20=25', '16+20=36'),\n  @('55+6=61', '20+34=54'),\n  @('1+39=40', '83-34=49'),\n  @('73+5=78', '24-15=9'),\n  @('88+4=92', '17+56=73'),\n  @('95-7=88', '70-31=39'),\n  @('14+42=56', '36+14=50'),\n  @('14+6=20', '25+73=98'),\n  @('32-23=9', '14+40=54'),\n  @('60-6=54', '38+46=84'),\n  @('63+6=69', '62+35=97'),\n  @('60-7=53', '55-37=18'),\n  @('42+22=64', '49+10=59'),\n  @('43-18=25', '28+26=54'),\n  @('34+3=37', '0+73=73'),\n  @('56+10=66', '69-57=12'),\n  @('42+19=61', '41+51=92'),\n  @('56-24=32', '19+76=95'),\n  @('20-17=3', '49+0=49'),\n  @('77+12=89', '25+64=89'),\n  @('1+59=60', '99-54=45'),\n  @('19+47=66', '36-30=6'),\n  @('83-73=10', '28+44=72'),\n  @('35-5=30', '88-6=82'),\n  @('19+5=24', '96-33=63'),\n  @('62-17=45', '34-4=30'),\n  @('40-8=32', '3+66=69'),\n  @('33+57=90', '93-79=14'),\n  @('21-5=16', '93-55=38'),\n  @('61-52=9', '37+0=37'),\n  @('91-39=52', '53-41=12'),\n  @('67+8=75', '15+58=73'),\n  @('44+45=89', '86-0=86'),\n  @('18+14=32', '67-22=45'),\n  @('21+46=67', '46+27=73'),\n  @('20+25=45', '18+67=85'),\n  @('37-21=16', '66-39=27'),\n  @('45-44=1', '85-13=72'),\n  @('15-3=12', '32+58=90'),\n  @('28-27=1', '15+65=80'),\n  @('11+34=45', '11+51=62'),\n  @('69-42=27', '20-9=11'),\n  @('15+59=74', '62-32=30'),\n  @('50-3=47', '24+44=68'),\n  @('37+23=60', '75-46=29'),\n  @('13+8=21', '73-25=48'),\n  @('44+31=75', '47-0=47'),\n  @('33+12=45', '16+43=59'),\n  @('58-40=18', '36+18=54'),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1          # wdFindContinue\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n\n  # wdReplaceAll = 2 : replace every (literal, unique) occurrence of $oldText\n  $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
